$wb = $excel.ActiveWorkbook

# --- Sheet1 updates ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 2: STATUS OK -> NO ISSUE
$ws1.Range("E2").Value = "NO ISSUE"

# Row 3: STATUS OK -> NO ISSUE, COMMENT tweak
$ws1.Range("E3").Value = "NO ISSUE"
$ws1.Range("F3").Value = "I think it's fine"

# Row 4: STATUS ERROR -> ISSUE, COMMENT tweak
$ws1.Range("E4").Value = "ISSUE"
$ws1.Range("F4").Value = "Missing article"

# --- Sheet2 updates ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 3: STATUS OK -> NO ISSUE, COMMENT tweak
$ws2.Range("E3").Value = "NO ISSUE"
$ws2.Range("F3").Value = "Shop verified"
